$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7050.6665
$ws.Range("J40").Value = 6860.8
$ws.Range("L40").Value = 6860.8
$ws.Range("N40").Value = -7210.8

$ws.Range("H51").Value = 3987.5

$ws.Range("H96").Value = 3421.7646
$ws.Range("I96").Value = 2383.8125
$ws.Range("K96").Value = 7151.4375
$ws.Range("M96").Value = -5778.4375

$ws.Range("H125").Value = 204247.6
$ws.Range("I125").Value = 335746
$ws.Range("J125").Value = 7000
$ws.Range("K125").Value = 3021714
$ws.Range("L125").Value = 63000
$ws.Range("M125").Value = -3019254
$ws.Range("N125").Value = -67920

$ws.Range("H137").Value = 280936.94
$ws.Range("I137").Value = 1856.7307
$ws.Range("J137").Value = 560017.1
$ws.Range("K137").Value = 5570.1921
$ws.Range("L137").Value = 1680051.3
$ws.Range("M137").Value = -3020.1921
$ws.Range("N137").Value = -1685151.3

$ws.Range("H138").Value = 1589.258
$ws.Range("I138").Value = 1249.45
$ws.Range("J138").Value = 1751.0714
$ws.Range("K138").Value = 3748.35
$ws.Range("L138").Value = 5253.2142
$ws.Range("M138").Value = 1391.65
$ws.Range("N138").Value = -15533.2142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9989.361000000001
$ws.Range("I32").Value = 5053.0527
$ws.Range("K32").Value = 5053.0527
$ws.Range("M32").Value = -4766.0527

$ws.Range("H45").Value = 11366042
$ws.Range("I45").Value = 2682.75
$ws.Range("K45").Value = 2682.75
$ws.Range("M45").Value = -2305.75

$ws.Range("H55").Value = 28999
$ws.Range("I55").Value = 28999
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 28999
$ws.Range("L55").ClearContents()
$ws.Range("N55").Value = 0
$ws.Range("M55").Value = -28684

$ws.Range("H74").Value = 2101.5386
$ws.Range("I74").Value = 1350.5
$ws.Range("J74").Value = 3303.2
$ws.Range("K74").Value = 1350.5
$ws.Range("L74").Value = 3303.2
$ws.Range("M74").Value = -476.5
$ws.Range("N74").Value = -5051.2

$ws.Range("H77").Value = 2101.5386
$ws.Range("I77").Value = 1350.5
$ws.Range("J77").Value = 3303.2
$ws.Range("K77").Value = 6752.5
$ws.Range("L77").Value = 16516
$ws.Range("M77").Value = -2384.5
$ws.Range("N77").Value = -25252

$ws.Range("H110").Value = 1171.5714
$ws.Range("I110").Value = 1171.5714
$ws.Range("K110").Value = 1171.5714
$ws.Range("M110").Value = 873.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 74434.64
$ws.Range("I105").Value = 145279.58
$ws.Range("K105").Value = 145279.58
$ws.Range("M105").Value = -143532.58

$ws.Range("H107").Value = 1147.4762
$ws.Range("I107").Value = 894.05554
$ws.Range("K107").Value = 894.05554
$ws.Range("M107").Value = 1025.94446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2274.3215
$ws.Range("I31").Value = 1691.3684
$ws.Range("J31").Value = 3505
$ws.Range("K31").Value = 1691.3684
$ws.Range("L31").Value = 3505
$ws.Range("M31").Value = -1396.3684
$ws.Range("N31").Value = -4095

$ws.Range("H34").Value = 2274.3215
$ws.Range("I34").Value = 1691.3684
$ws.Range("J34").Value = 3505
$ws.Range("K34").Value = 1691.3684
$ws.Range("L34").Value = 3505
$ws.Range("M34").Value = -1489.3684
$ws.Range("N34").Value = -3909

$ws.Range("H60").Value = 3871.375
$ws.Range("I60").Value = 3871.375
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 3871.375
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -3360.375

$ws.Range("H99").Value = 11113795
$ws.Range("J99").Value = 3166.6667
$ws.Range("L99").Value = 3166.6667
$ws.Range("N99").Value = -6162.6667

$ws.Range("H107").Value = 1185.7
$ws.Range("I107").Value = 1204.8889
$ws.Range("K107").Value = 1204.8889
$ws.Range("M107").Value = 715.1111000000001

$ws.Range("H122").Value = 2821.6875
$ws.Range("I122").Value = 2360.8
$ws.Range("K122").Value = 7082.400000000001
$ws.Range("M122").Value = -4632.400000000001

$ws.Range("H126").Value = 11113795
$ws.Range("J126").Value = 3166.6667
$ws.Range("L126").Value = 9500.000100000001
$ws.Range("N126").Value = -14440.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 385811
$ws.Range("I128").Value = 385811
$ws.Range("K128").Value = 1157433
$ws.Range("M128").Value = -1152453

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 3140.3076
$ws.Range("I22").Value = 3225
$ws.Range("J22").Value = 3041.5
$ws.Range("K22").Value = 3225
$ws.Range("L22").Value = 3041.5
$ws.Range("M22").Value = -2696
$ws.Range("N22").Value = -4099.5

$ws.Range("H113").Value = 2240
$ws.Range("I113").Value = 1733.3334
$ws.Range("K113").Value = 1733.3334
$ws.Range("M113").Value = 436.6666

$ws.Range("H122").Value = 373598.4
$ws.Range("I122").Value = 529065.7
$ws.Range("K122").Value = 1587197.1
$ws.Range("M122").Value = -1584747.1

$ws.Range("H133").Value = 64796.8
$ws.Range("J133").Value = 70996
$ws.Range("L133").Value = 70996
$ws.Range("N133").Value = -81116

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 5259.636
$ws.Range("J55").Value = 18051.5
$ws.Range("L55").Value = 18051.5
$ws.Range("N55").Value = -18397.5

$ws.Range("H82").Value = 2958.8
$ws.Range("I82").Value = 3011
$ws.Range("J82").Value = 2750
$ws.Range("K82").Value = 3011
$ws.Range("L82").Value = 2750
$ws.Range("M82").Value = -2650
$ws.Range("N82").Value = -3472

$ws.Range("H85").Value = 2958.8
$ws.Range("I85").Value = 3011
$ws.Range("J85").Value = 2750
$ws.Range("K85").Value = 3011
$ws.Range("L85").Value = 2750
$ws.Range("M85").Value = -1763
$ws.Range("N85").Value = -5246

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 3000
$ws.Range("I39").Value = 3000
$ws.Range("K39").Value = 3000
$ws.Range("M39").Value = -2587

$ws.Range("H107").Value = 1515
$ws.Range("J107").Value = 3294.1428
$ws.Range("L107").Value = 9882.428400000001
$ws.Range("N107").Value = -13722.4284

$ws.Range("H113").Value = 1848.8
$ws.Range("I113").Value = 1835.909
$ws.Range("J113").Value = 1864.5555
$ws.Range("K113").Value = 5507.727000000001
$ws.Range("L113").Value = 5593.666499999999
$ws.Range("M113").Value = -3337.727000000001
$ws.Range("N113").Value = -9933.666499999999

$ws.Range("H126").Value = 17570.857
$ws.Range("I126").Value = 26249
$ws.Range("K126").Value = 78747
$ws.Range("M126").Value = -76277
